# Apply crypto price/volume updates (and an Algorand/Elrond row swap)
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.450.43"
$ws.Range("E2").Value = "'  -0.94%  "
$ws.Range("D3").Value = "'1.915.13"
$ws.Range("E3").Value = "'  +1.72%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("D5").Value = "'241.14"
$ws.Range("E5").Value = "'  +1.27%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.12%  "
$ws.Range("D7").Value = "'0.4714"
$ws.Range("E7").Value = "'  -0.78%  "
$ws.Range("E8").Value = "'  +1.17%  "
$ws.Range("D9").Value = "'0.06799"
$ws.Range("D10").Value = "'106.35"
$ws.Range("E10").Value = "'  +11.86%  "
$ws.Range("D11").Value = "'18.39"
$ws.Range("E11").Value = "'  -1.12%  "
$ws.Range("D12").Value = "'1.904.67"
$ws.Range("E12").Value = "'  +1.02%  "
$ws.Range("D13").Value = "'0.07700"
$ws.Range("E13").Value = "'  +1.84%  "
$ws.Range("D14").Value = "'5.213"
$ws.Range("E14").Value = "'  +2.77%  "
$ws.Range("D15").Value = "'0.6577"
$ws.Range("E15").Value = "'  +1.36%  "
$ws.Range("D16").Value = "'289.93"
$ws.Range("E16").Value = "'  -3.84%  "
$ws.Range("D17").Value = "'30.460.89"
$ws.Range("E17").Value = "'  -0.86%  "
$ws.Range("D18").Value = "'0.000007639"
$ws.Range("E18").Value = "'  +1.31%  "
$ws.Range("E19").Value = "'  +0.04%  "
$ws.Range("E20").Value = "'  -0.88%  "
$ws.Range("D21").Value = "'2.154.27"
$ws.Range("E21").Value = "'  +1.40%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "'  +0.11%  "
$ws.Range("D23").Value = "'5.211"
$ws.Range("E23").Value = "'  +1.55%  "
$ws.Range("D24").Value = "'6.213"
$ws.Range("E24").Value = "'  +1.15%  "
$ws.Range("D25").Value = "'9.319"
$ws.Range("E25").Value = "'  +1.00%  "
$ws.Range("D26").Value = "'168.13"
$ws.Range("E26").Value = "'  -0.60%  "
$ws.Range("D27").Value = "'21.52"
$ws.Range("E27").Value = "'  +9.34%  "
$ws.Range("D28").Value = "'2.084"
$ws.Range("E28").Value = "'  +7.36%  "
$ws.Range("D29").Value = "'0.1067"
$ws.Range("E29").Value = "'  +0.48%  "
$ws.Range("E30").Value = "'  +1.33%  "
$ws.Range("D31").Value = "'4.159"
$ws.Range("E31").Value = "'  +0.15%  "
$ws.Range("D32").Value = "'3.983"
$ws.Range("E32").Value = "'  +1.08%  "
$ws.Range("D33").Value = "'0.05059"
$ws.Range("E33").Value = "'  -0.01%  "
$ws.Range("D34").Value = "'0.7440"
$ws.Range("E34").Value = "'  +3.26%  "
$ws.Range("D35").Value = "'1.152"
$ws.Range("E35").Value = "'  -1.27%  "
$ws.Range("D36").Value = "'0.02096"
$ws.Range("E36").Value = "'  +9.41%  "
$ws.Range("D37").Value = "'2.740"
$ws.Range("E37").Value = "'  +0.80%  "
$ws.Range("D38").Value = "'2.672"
$ws.Range("E38").Value = "'  -1.37%  "
$ws.Range("D39").Value = "'2.055"
$ws.Range("E39").Value = "'  +0.21%  "
$ws.Range("D40").Value = "'109.55"
$ws.Range("E40").Value = "'  +2.07%  "
$ws.Range("D41").Value = "'0.8702"
$ws.Range("E41").Value = "'  -3.03%  "
$ws.Range("D42").Value = "'5.882"
$ws.Range("E42").Value = "'  +5.09%  "
$ws.Range("D43").Value = "'0.4271"
$ws.Range("E43").Value = "'  +2.13%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "'  +0.09%  "
$ws.Range("D45").Value = "'67.68"
$ws.Range("E45").Value = "'  +4.43%  "
$ws.Range("D46").Value = "'50.55"
$ws.Range("E46").Value = "'  +18.73%  "
$ws.Range("D47").Value = "'7.178"
$ws.Range("E47").Value = "'  -1.67%  "
$ws.Range("D48").Value = "'9.268"
$ws.Range("E48").Value = "'  +3.50%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'34.95"
$ws.Range("E49").Value = "'  +1.15%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1213"
$ws.Range("E50").Value = "'  -0.39%  "
$ws.Range("D51").Value = "'0.3931"
$ws.Range("E51").Value = "'  +3.28%  "
